# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 9
$ws1.Range("F3").Value  = 171
$ws1.Range("F7").Value  = 1663
$ws1.Range("F10").Value = 26
$ws1.Range("F11").Value = 1585
$ws1.Range("F13").Value = 59
$ws1.Range("F14").Value = 394
$ws1.Range("F15").Value = 259
$ws1.Range("F18").Value = 20
$ws1.Range("F19").Value = 29
$ws1.Range("F20").Value = 49
$ws1.Range("F21").Value = 64
$ws1.Range("F22").Value = 283
$ws1.Range("F24").Value = 222
$ws1.Range("F25").Value = 221

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 9
$ws4.Range("F3").Value  = 171
$ws4.Range("F7").Value  = 1663
$ws4.Range("F11").Value = 26
$ws4.Range("F12").Value = 1585
$ws4.Range("F14").Value = 59
$ws4.Range("F15").Value = 394
$ws4.Range("F16").Value = 259
$ws4.Range("F19").Value = 20
$ws4.Range("F20").Value = 29
$ws4.Range("F21").Value = 49
$ws4.Range("F22").Value = 64
$ws4.Range("F23").Value = 283
$ws4.Range("F25").Value = 222
$ws4.Range("F26").Value = 221
